$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# String index lookup (shared strings are 0-based: 20=ECs, 21=FAPs, 22=MuSCs, 23=Vcan, 24=Selp)
$ECs = "ECs"
$FAPs = "FAPs"
$MuSCs = "MuSCs"
$Vcan = "Vcan"
$Selp = "Selp"

# Row 2: ECs -> ECs
$ws.Range("A2").Value = $ECs
$ws.Range("B2").Value = $Vcan
$ws.Range("C2").Value = $Selp
$ws.Range("D2").Value = $ECs
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2364713333333333
$ws.Range("H2").Value = 0.709414
$ws.Range("I2").Value = 0.002249544876489787
$ws.Range("J2").Value = 0.002249544876489787
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.866432
$ws.Range("N2").Value = 8.599295999999999
$ws.Range("O2").Value = 0.9456981836489474
$ws.Range("P2").Value = 0.9456981836489475
$ws.Range("Q2").Value = 0.6778289969493333
$ws.Range("R2").Value = 6.100460972543999
$ws.Range("S2").Value = 0.002127390503733188
$ws.Range("T2").Value = 0.002127390503733188

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = $ECs
$ws.Range("B3").Value = $Vcan
$ws.Range("C3").Value = $Selp
$ws.Range("D3").Value = $FAPs
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2364713333333333
$ws.Range("H3").Value = 0.709414
$ws.Range("I3").Value = 0.002249544876489787
$ws.Range("J3").Value = 0.002249544876489787
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.16459
$ws.Range("N3").Value = 0.49377
$ws.Range("O3").Value = 0.05430181635105255
$ws.Range("P3").Value = 0.05430181635105256
$ws.Range("Q3").Value = 0.03892081675333333
$ws.Range("R3").Value = 0.35028735078
$ws.Range("S3").Value = 0.0001221543727565996
$ws.Range("T3").Value = 0.0001221543727565996

# Row 4: FAPs -> ECs
$ws.Range("A4").Value = $FAPs
$ws.Range("B4").Value = $Vcan
$ws.Range("C4").Value = $Selp
$ws.Range("D4").Value = $ECs
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 103.676216
$ws.Range("H4").Value = 311.028648
$ws.Range("I4").Value = 0.9862688099613843
$ws.Range("J4").Value = 0.9862688099613843
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.866432
$ws.Range("N4").Value = 8.599295999999999
$ws.Range("O4").Value = 0.9456981836489474
$ws.Range("P4").Value = 0.9456981836489475
$ws.Range("Q4").Value = 297.180823181312
$ws.Range("R4").Value = 2674.627408631808
$ws.Range("S4").Value = 0.9327126221700901
$ws.Range("T4").Value = 0.9327126221700902

# Row 5: FAPs -> FAPs
$ws.Range("A5").Value = $FAPs
$ws.Range("B5").Value = $Vcan
$ws.Range("C5").Value = $Selp
$ws.Range("D5").Value = $FAPs
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 103.676216
$ws.Range("H5").Value = 311.028648
$ws.Range("I5").Value = 0.9862688099613843
$ws.Range("J5").Value = 0.9862688099613843
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.16459
$ws.Range("N5").Value = 0.49377
$ws.Range("O5").Value = 0.05430181635105255
$ws.Range("P5").Value = 0.05430181635105256
$ws.Range("Q5").Value = 17.06406839144
$ws.Range("R5").Value = 153.57661552296
$ws.Range("S5").Value = 0.05355618779129424
$ws.Range("T5").Value = 0.05355618779129425

# Row 6: MuSCs -> ECs
$ws.Range("A6").Value = $MuSCs
$ws.Range("B6").Value = $Vcan
$ws.Range("C6").Value = $Selp
$ws.Range("D6").Value = $ECs
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.206946333333333
$ws.Range("H6").Value = 3.620839
$ws.Range("I6").Value = 0.01148164516212593
$ws.Range("J6").Value = 0.01148164516212593
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.866432
$ws.Range("N6").Value = 8.599295999999999
$ws.Range("O6").Value = 0.9456981836489474
$ws.Range("P6").Value = 0.9456981836489475
$ws.Range("Q6").Value = 3.459629592149333
$ws.Range("R6").Value = 31.136666329344
$ws.Range("S6").Value = 0.01085817097512422
$ws.Range("T6").Value = 0.01085817097512422

# Row 7: MuSCs -> FAPs
$ws.Range("A7").Value = $MuSCs
$ws.Range("B7").Value = $Vcan
$ws.Range("C7").Value = $Selp
$ws.Range("D7").Value = $FAPs
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.206946333333333
$ws.Range("H7").Value = 3.620839
$ws.Range("I7").Value = 0.01148164516212593
$ws.Range("J7").Value = 0.01148164516212593
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.16459
$ws.Range("N7").Value = 0.49377
$ws.Range("O7").Value = 0.05430181635105255
$ws.Range("P7").Value = 0.05430181635105256
$ws.Range("Q7").Value = 0.1986512970033333
$ws.Range("R7").Value = 1.78786167303
$ws.Range("S7").Value = 0.0006234741870017133
$ws.Range("T7").Value = 0.0006234741870017134

# Remove rows 8-10 entirely (previously had data, now removed)
$ws.Range("A8:T10").Delete()
